$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Robo2","Robo2","ECs",3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,0.08193410090844443,0.7374069081759999,0.06103278096779807,0.06103278096779807),
    @("ECs","Robo2","Robo2","FAPs",3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,0.2490610648915555,2.241549584023999,0.185525797593346,0.185525797593346),
    @("ECs","Robo2","Robo2","sCs",3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,0.0006572101013333332,0.005914890912,0.0004895563595592917,0.0004895563595592917),
    @("FAPs","Robo2","Robo2","ECs",3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,0.2490610648915555,2.241549584023999,0.185525797593346,0.185525797593346),
    @("FAPs","Robo2","Robo2","FAPs",3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,0.7570890918084442,6.813801826275999,0.5639563039214562,0.5639563039214562),
    @("FAPs","Robo2","Robo2","sCs",3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,0.001997769498666666,0.017979925488,0.001488140186862097,0.001488140186862097),
    @("sCs","Robo2","Robo2","ECs",2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,3,1,0.2862413333333333,0.8587239999999999,0.2470481349207034,0.2470481349207034,0.0006572101013333332,0.005914890912,0.0004895563595592917,0.0004895563595592917),
    @("sCs","Robo2","Robo2","FAPs",2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,3,1,0.8701086666666665,2.610326,0.7509702417016644,0.7509702417016644,0.001997769498666666,0.017979925488,0.001488140186862097,0.001488140186862097),
    @("sCs","Robo2","Robo2","sCs",2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,2,0.6666666666666666,0.002296,0.006888,0.001981623377632167,0.001981623377632167,0.000005271615999999999,0.000047444544,0.00000392683121077832,0.00000392683121077832)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $row[$j]
    }
}
